$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need their number format
# forced to Text first, otherwise Excel would auto-convert the text into a
# numeric value. The style is reset back to Normal afterwards so the cell
# keeps its original (default) styling, matching the source workbook.
$numericLookingCells = @('D5', 'D8', 'D9', 'D13', 'D17', 'D18', 'D20', 'D21', 'D24', 'D26', 'D27', 'D32', 'D34', 'D39', 'D40', 'D41', 'D43', 'D44', 'D46', 'D49')
foreach ($cell in $numericLookingCells) {
    $ws.Range($cell).NumberFormat = '@'
}

$ws.Range('D2').Value = '34.469.73'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '1.805.63'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('D5').Value = '228.48'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +4.63%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '34.93'
$ws.Range('E8').Value = '  +6.12%  '
$ws.Range('D9').Value = '0.301'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '2.065.80'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = '11.26'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('D14').Value = '1.807.06'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('D16').Value = '34.454.66'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '4.37'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = '69.14'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '0.0₃0798'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').Value = '246.00'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').Value = '11.52'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '173.73'
$ws.Range('E24').Value = '  +4.48%  '
$ws.Range('E25').Value = '  +2.24%  '
$ws.Range('D26').Value = '7.71'
$ws.Range('E26').Value = '  +5.60%  '
$ws.Range('D27').Value = '16.81'
$ws.Range('E27').Value = '  +1.49%  '
$ws.Range('E28').Value = '  +2.61%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('D32').Value = '3.84'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').Value = '1.397.49'
$ws.Range('E35').Value = '  -2.14%  '
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = '0.0191'
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('D40').Value = '83.41'
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('D41').Value = '2.84'
$ws.Range('E41').Value = '  +3.12%  '
$ws.Range('E42').Value = '  +1.67%  '
$ws.Range('D43').Value = '2.39'
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').Value = '13.56'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('E45').Value = '  +3.34%  '
$ws.Range('D46').Value = '0.0510'
$ws.Range('E46').Value = '  -3.38%  '
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('D48').Value = '1.965.73'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').Value = '105.03'
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('E51').Value = '  +0.50%  '

foreach ($cell in $numericLookingCells) {
    $ws.Range($cell).Style = 'Normal'
}
